# System/system_management.xlsx
# "asignamos WoS en base a clase social"
#
# Reassign the "archetypes" labels (column A) to their correct rows and
# add a new "stats" column (F) that records the social-class-based
# classification (class / synpop / trans) used to assign WoS activities.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: fix archetype labels on rows 2-6 -------------------------
# (the B/C/D/E data per row stays put; only the archetype name moves)
$ws.Range("A2").Value = "building"
$ws.Range("A3").Value = "citizen"
$ws.Range("A4").Value = "distribution"
$ws.Range("A5").Value = "family"
$ws.Range("A6").Value = "transport"

# --- New column F: "stats" ------------------------------------------------
$ws.Range("F1").Value = "stats"
$ws.Range("F2").Value = "synpop"
$ws.Range("F3").Value = "trans"
$ws.Range("F4").Value = "class"

# Sort the new stats column (data rows only, header excluded) alphabetically
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("F1:F4"))
$sortObj.SetRange($ws.Range("F1:F4"))
$sortObj.Header = 1
$sortObj.Apply()

# --- Header row formatting: bold -----------------------------------------
$ws.Range("A1:F1").Font.Bold = $true

# --- Column widths (approx. autofit to the new content) ------------------
$ws.Columns.Item(1).ColumnWidth = 9.6
$ws.Columns.Item(2).ColumnWidth = 9.28
$ws.Columns.Item(3).ColumnWidth = 11.83
$ws.Columns.Item(4).ColumnWidth = 2.94
$ws.Columns.Item(5).ColumnWidth = 12.05
$ws.Columns.Item(6).ColumnWidth = 5.94

# --- View: zoom to 210% and select the full A:F columns ------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 210
$ws.Range("A1:F1048576").Select()
